$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column (15-nov) before column DT ---
$ws1 = $wb.Worksheets.Item("Prix Spot")
$ws1.Range("DT1").EntireColumn.Insert()
$ws1.Range("DT1").Value = "15-nov"
$ws1.Range("DT2:DT25").Value = "-"

# --- Sheet "CO2": append new daily price row ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A152").NumberFormat = "@"
$ws3.Range("A152").Value = "2025-11-13"
$ws3.Range("A152").Style = "Normal"
$ws3.Range("B152").Value = 81.02
